$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2021-12-25"

# Update the label in A14
$ws.Range("A14").Value = "December (through 12-25)"

# Row 14 updates
$ws.Range("C14").Value = 32
$ws.Range("D14").Value = 0.1111
$ws.Range("F14").Value = 75
$ws.Range("G14").Value = 0.0854
$ws.Range("H14").Value = 11
$ws.Range("I14").Value = 86
$ws.Range("J14").Value = 0.1134
$ws.Range("L14").Value = 52
$ws.Range("M14").Value = 0.0877
$ws.Range("O14").Value = 48
$ws.Range("P14").Value = 0.0769
$ws.Range("R14").Value = 110
$ws.Range("S14").Value = 0.0678
$ws.Range("U14").Value = 156
$ws.Range("V14").Value = 0.0127

# Row 15 updates (Totals)
$ws.Range("C15").Value = 290
$ws.Range("D15").Value = 0.1131
$ws.Range("F15").Value = 579
$ws.Range("G15").Value = 0.1023
$ws.Range("H15").Value = 74
$ws.Range("I15").Value = 844
$ws.Range("J15").Value = 0.0806
$ws.Range("L15").Value = 660
$ws.Range("M15").Value = 0.1069
$ws.Range("O15").Value = 528
$ws.Range("P15").Value = 0.099
$ws.Range("R15").Value = 1310
$ws.Range("S15").Value = 0.0521
$ws.Range("U15").Value = 1699
$ws.Range("V15").Value = 0.0566
